$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    ,@(39400,2007,4.930115226412335,2008,1.402909115503936)
    ,@(39583,2008,2.319409367208825,2009,5.080273296954374)
    ,@(39765,2008,1.457587285166628,2009,2.807231216534278)
    ,@(39948,2009,-3.942037578692481,2010,-1.648748515828502)
    ,@(40130,2009,-0.9140166223623569,2010,1.821983295885121)
    ,@(40310,2010,-2.839753013810475,2011,-1.632723506456935)
    ,@(40494,2010,2.585942866987878,2011,2.722861752007866)
    ,@(40676,2011,4.960109259035406,2012,4.506881698240095)
    ,@(40862,2011,4.253963781362402,2012,2.878414118480799)
    ,@(41044,2012,2.134646894829806,2013,2.531943146540772)
    ,@(41228,2012,1.752870900283909,2013,3.144721336271927)
    ,@(41409,2013,-2.700325749999488,2014,-0.3858735870725938)
    ,@(41592,2013,-1.479696720105139,2014,2.238623952069552)
    ,@(41774,2014,5.469647210234996,2015,3.061326532789543)
    ,@(41957,2014,3.900127535411246,2015,-0.6155071485167585)
    ,@(42137,2015,0.950153436409007,2016,2.074800935750787)
    ,@(42321,2015,0.03947433952959933,2016,1.182212550358064)
    ,@(42503,2016,3.458696398997052,2017,2.610227683091337)
    ,@(42689,2016,2.192778679161944,2017,-0.5835597102573198)
    ,@(42867,2017,2.772413308959698,2018,2.755099409670958)
    ,@(43053,2017,3.40836448860673,2018,2.34299484087257)
    ,@(43145,2018,1.376958470962375,2019,2.554746761760573)
    ,@(43235,2018,3.145819842658448,2019,3.875106770584158)
    ,@(43326,2018,2.594480907596441,2019,2.978856953395592)
    ,@(43418,2018,2.799070570134488,2019,3.056075254339996)
    ,@(43510,2019,3.665688413913726,2020,3.126710782028086)
    ,@(43600,2019,4.520465362328063,2020,4.124307769579505)
    ,@(43691,2019,4.143226503463837,2020,2.872814227731846)
    ,@(43783,2019,4.195393191694419,2020,2.652948310315506)
    ,@(43875,2020,2.017172445889925,2021,2.549404806734712)
    ,@(43966,2020,4.479055418855871,2021,4.888255652935936)
    ,@(44068,2020,2.133862376612439,2021,0.273062749117492)
    ,@(44159,2020,1.666553973046048,2021,-0.9999522486825452)
    ,@(44251,2021,-2.991914460776179,2022,-0.7897024150736676)
    ,@(44341,2021,2.922623512367206,2022,2.761298099516418)
    ,@(44432,2021,1.773820722495745,2022,1.707328096006822)
    ,@(44525,2021,1.879266440112803,2022,-0.1343977949472275)
    ,@(44617,2022,-0.05874878927536464,2023,1.647170385765651)
    ,@(44706,2022,-0.2939924376199055,2023,1.109922826784815)
    ,@(44798,2022,-2.404913754290983,2023,-2.467583266971451)
    ,@(44890,2022,-2.620683231370946,2023,-2.83913279674276)
    ,@(44981,2023,-3.938880009788048,2024,-0.6223811084860253)
    ,@(45071,2023,-2.429883624035745,2024,-0.8739913853863412)
    ,@(45163,2023,-2.901570548279864,2024,-1.444783760697688)
    ,@(45254,2023,-3.036556262700274,2024,-1.803491225663911)
    ,@(45345,2024,-2.786207191231715,2025,-2.423357831382655)
    ,@(45436,2024,-1.060428249734879,2025,-1.285332554730334)
    ,@(45534,2024,-2.859191689251428,2025,-3.927752965551978)
    ,@(45618,2024,-2.953443685011514,2025,-2.75492543068685)
    ,@(45713,2025,-0.8799002368073738,2026,-2.26429946512674)
    ,@(45800,2025,0.1432745330888219,2026,-0.484131235569496)
    ,@(45891,2025,-1.131442475565558,2026,-2.567308877837815)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $data[$i]
    $r = $i + 2
    for ($j = 0; $j -lt $row.Count; $j++) {
        $c = $j + 1
        $val = $row[$j]
        if ($null -ne $val) {
            $ws.Cells.Item($r, $c).Value = $val
        }
    }
}

# Row 53 is newly added; give its date cell (column A) the same
# formatting (date number format / font / border / alignment) as the
# preceding rows' date cells by copying the format down.
$ws.Cells.Item(52, 1).Copy() | Out-Null
$ws.Cells.Item(53, 1).PasteSpecial(-4122) | Out-Null
